# Generate Report for Handoff
#
# The ec14b915-... file has moved from "Handed back: in sync with en-US"
# to "Ready for handoff" status, with refreshed handoff timestamps and a
# new error detail explaining the handback file is stale. Reflected on
# the Overview sheet as well as both the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2cf0ca8398693ca63c348dbdd9c47766f33e4b2/e2e/ec14b915-a61b-4c66-8d3f-aeef0ee0add2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee2f0665099a381f9fa895c2b05dad5694e6ff39/e2e/ec14b915-a61b-4c66-8d3f-aeef0ee0add2.md."

# --- Overview sheet: row 3 is the ec14b915 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 22:48:10"

# --- zh-cn sheet: row 3 is the ec14b915 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-28 22:48:06"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39 + 1/6

# --- de-de sheet: row 3 is the ec14b915 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-28 22:48:10"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39 + 1/6
